$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-23 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-24 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("741×7=5187", $true, $false, $false, $false, $false, $true, 1, $false, "444×9=3996", 2) | Out-Null
$d.Content.Find.Execute("646×5=3230", $true, $false, $false, $false, $false, $true, 1, $false, "659×8=5272", 2) | Out-Null
$d.Content.Find.Execute("112×4=448", $true, $false, $false, $false, $false, $true, 1, $false, "916×8=7328", 2) | Out-Null
$d.Content.Find.Execute("980×9=8820", $true, $false, $false, $false, $false, $true, 1, $false, "649×3=1947", 2) | Out-Null
$d.Content.Find.Execute("203×2=406", $true, $false, $false, $false, $false, $true, 1, $false, "223×2=446", 2) | Out-Null
$d.Content.Find.Execute("503×8=4024", $true, $false, $false, $false, $false, $true, 1, $false, "666×2=1332", 2) | Out-Null
$d.Content.Find.Execute("646×3=1938", $true, $false, $false, $false, $false, $true, 1, $false, "729×6=4374", 2) | Out-Null
$d.Content.Find.Execute("671×5=3355", $true, $false, $false, $false, $false, $true, 1, $false, "988×2=1976", 2) | Out-Null
$d.Content.Find.Execute("558×8=4464", $true, $false, $false, $false, $false, $true, 1, $false, "578×2=1156", 2) | Out-Null
$d.Content.Find.Execute("883×3=2649", $true, $false, $false, $false, $false, $true, 1, $false, "138×8=1104", 2) | Out-Null
$d.Content.Find.Execute("447×4=1788", $true, $false, $false, $false, $false, $true, 1, $false, "293×3=879", 2) | Out-Null
$d.Content.Find.Execute("255×3=765", $true, $false, $false, $false, $false, $true, 1, $false, "553×3=1659", 2) | Out-Null
$d.Content.Find.Execute("152×3=456", $true, $false, $false, $false, $false, $true, 1, $false, "853×4=3412", 2) | Out-Null
$d.Content.Find.Execute("332×3=996", $true, $false, $false, $false, $false, $true, 1, $false, "876×4=3504", 2) | Out-Null
$d.Content.Find.Execute("876×8=7008", $true, $false, $false, $false, $false, $true, 1, $false, "596×6=3576", 2) | Out-Null
$d.Content.Find.Execute("980×2=1960", $true, $false, $false, $false, $false, $true, 1, $false, "379×5=1895", 2) | Out-Null
$d.Content.Find.Execute("212×6=1272", $true, $false, $false, $false, $false, $true, 1, $false, "310×9=2790", 2) | Out-Null
$d.Content.Find.Execute("985×5=4925", $true, $false, $false, $false, $false, $true, 1, $false, "858×7=6006", 2) | Out-Null
$d.Content.Find.Execute("459×2=918", $true, $false, $false, $false, $false, $true, 1, $false, "620×5=3100", 2) | Out-Null
$d.Content.Find.Execute("774×9=6966", $true, $false, $false, $false, $false, $true, 1, $false, "764×7=5348", 2) | Out-Null
$d.Content.Find.Execute("747×9=6723", $true, $false, $false, $false, $false, $true, 1, $false, "750×5=3750", 2) | Out-Null
$d.Content.Find.Execute("640×9=5760", $true, $false, $false, $false, $false, $true, 1, $false, "592×3=1776", 2) | Out-Null
$d.Content.Find.Execute("998×5=4990", $true, $false, $false, $false, $false, $true, 1, $false, "950×4=3800", 2) | Out-Null
$d.Content.Find.Execute("515×6=3090", $true, $false, $false, $false, $false, $true, 1, $false, "611×6=3666", 2) | Out-Null
$d.Content.Find.Execute("543×8=4344", $true, $false, $false, $false, $false, $true, 1, $false, "562×8=4496", 2) | Out-Null
